$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Shared text "2016-08-20 00:14:33" -> "2016-08-20 00:15:19"
# appears as Overview "Latest HO Xliff Generate Date" (G2, G5)
# and de-de "Correspond Handoff Datetime" (H2, H5)
$wsOverview.Range("G2").Value = "2016-08-20 00:15:19"
$wsOverview.Range("G5").Value = "2016-08-20 00:15:19"
$wsDeDe.Range("H2").Value = "2016-08-20 00:15:19"
$wsDeDe.Range("H5").Value = "2016-08-20 00:15:19"

# Shared text "ht" -> "mt" (Status column E) on both zh-cn and de-de
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn "Correspond Handoff Datetime" (H2, H5): "2016-08-20 00:14:29" -> "2016-08-20 00:15:15"
$wsZhCn.Range("H2").Value = "2016-08-20 00:15:15"
$wsZhCn.Range("H5").Value = "2016-08-20 00:15:15"

# zh-cn "Correspond Handback DateTime" (K2, K5): "2016-08-20 00:14:45" -> "2016-08-20 00:15:32"
$wsZhCn.Range("K2").Value = "2016-08-20 00:15:32"
$wsZhCn.Range("K5").Value = "2016-08-20 00:15:32"

# de-de "Correspond Handback DateTime" (K2, K5): "2016-08-20 00:14:51" -> "2016-08-20 00:15:39"
$wsDeDe.Range("K2").Value = "2016-08-20 00:15:39"
$wsDeDe.Range("K5").Value = "2016-08-20 00:15:39"
